# Delete the entire row 558 ("自動車メーカー" / car manufacturers post),
# which shifts all subsequent rows up by one and reduces the sheet's
# used range from A1:C655 to A1:C654.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(558).Delete()
